# Add a "version" column (and corresponding "version list" sheet) to the
# MALDI-IMS metadata template, per commit "Add version 1 everywhere".
#
# Summary of the edit:
#   1. A new worksheet "version list" is inserted right after "Export as TSV"
#      (i.e. as the 2nd sheet), containing the single value "1" in A1.
#   2. On the "Export as TSV" sheet, a new first column "version" is
#      inserted before the existing "donor_id" column, shifting every other
#      column one position to the right.
#   3. Every existing header-cell comment shifts one column to the right
#      along with its column, and a new comment is added on the new A1
#      ("version") cell.
#   4. A new data validation (list, sourced from 'version list'!$A$1:$A$1)
#      is added for column A; all pre-existing data validations keep
#      their original relative column but shift right by one letter.

function Col-Letter($n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - $rem) / 26)
    }
    return $s
}

$wb = $excel.ActiveWorkbook
$tsv = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------
# 1. New "version list" sheet, placed right after "Export as TSV".
# ---------------------------------------------------------------------
$verList = $wb.Worksheets.Add($null, $tsv)
$verList.Name = "version list"
$verList.Range("A1").NumberFormat = "@"
$verList.Range("A1").Value = "1"

# ---------------------------------------------------------------------
# 2. Insert the new "version" column at the front of the main sheet.
#    Columns.Insert() shifts cell values AND the sqref of existing
#    dataValidation entries one column to the right automatically, but
#    it does NOT move cell comments - those are handled separately below.
# ---------------------------------------------------------------------
$tsv.Columns.Item(1).Insert()

$tsv.Range("A1").Value = "version"
$tsv.Range("A1").Font.Bold = $true
$tsv.Range("A1").HorizontalAlignment = -4108
$tsv.Range("A1").WrapText = $true

# ---------------------------------------------------------------------
# 3. Shift every header comment one column to the right (process from
#    the last column back to the first so we never clobber a comment
#    before it has been read), then add the new comment for A1.
# ---------------------------------------------------------------------
for ($i = 30; $i -ge 1; $i--) {
    $srcCol = Col-Letter($i)
    $dstCol = Col-Letter($i + 1)
    $srcCell = $tsv.Range($srcCol + "1")
    $commentText = $srcCell.Comment.Text()
    $srcCell.Comment.Delete()
    $dstCell = $tsv.Range($dstCol + "1")
    $dstCell.AddComment($commentText)
}

$tsv.Range("A1").AddComment("Current version of metadata schema. Template provides the correct value.")

# ---------------------------------------------------------------------
# 4. Rebuild data validations in the same relative order as the source
#    diff: the new "version" validation first, followed by all the
#    pre-existing validations (already shifted one column right by the
#    earlier Columns.Insert() call).
# ---------------------------------------------------------------------
$existing = @(
    @{ Col = "J"; Type = "list"; Formula1 = "'assay_category list'!`$A`$1:`$A`$1"; Title = "Value must come from list"; Msg = "Value must be one of: mass_spectrometry_imaging." },
    @{ Col = "K"; Type = "list"; Formula1 = "'assay_type list'!`$A`$1:`$A`$1"; Title = "Value must come from list"; Msg = "Value must be one of: MALDI-IMS." },
    @{ Col = "L"; Type = "list"; Formula1 = "'analyte_class list'!`$A`$1:`$A`$3"; Title = "Value must come from list"; Msg = "Value must be one of: protein / metabolites / lipids." },
    @{ Col = "M"; Type = "list"; Formula1 = '"TRUE,FALSE"'; Title = "Not a boolean"; Msg = 'The values in this column must be "TRUE" or "FALSE".' },
    @{ Col = "P"; Type = "list"; Formula1 = "'ms_source list'!`$A`$1:`$A`$5"; Title = "Value must come from list"; Msg = "Value must be one of: MALDI / MALDI-2 / DESI / SIMS / nESI." },
    @{ Col = "Q"; Type = "list"; Formula1 = "'polarity list'!`$A`$1:`$A`$2"; Title = "Value must come from list"; Msg = "Value must be one of: negative ion mode / positive ion mode." },
    @{ Col = "R"; Type = "decimal"; Formula1 = "-1e+307"; Formula2 = "1e+307"; Title = "Not a number"; Msg = "The values in this column must be numbers." },
    @{ Col = "S"; Type = "decimal"; Formula1 = "-1e+307"; Formula2 = "1e+307"; Title = "Not a number"; Msg = "The values in this column must be numbers." },
    @{ Col = "T"; Type = "decimal"; Formula1 = "-1e+307"; Formula2 = "1e+307"; Title = "Not a number"; Msg = "The values in this column must be numbers." },
    @{ Col = "U"; Type = "list"; Formula1 = "'resolution_x_unit list'!`$A`$1:`$A`$2"; Title = "Value must come from list"; Msg = "Value must be one of: nm / um." },
    @{ Col = "V"; Type = "decimal"; Formula1 = "-1e+307"; Formula2 = "1e+307"; Title = "Not a number"; Msg = "The values in this column must be numbers." },
    @{ Col = "W"; Type = "list"; Formula1 = "'resolution_y_unit list'!`$A`$1:`$A`$2"; Title = "Value must come from list"; Msg = "Value must be one of: nm / um." }
)

$fullRange = $tsv.Range("A2:W1048576")
$fullRange.Validation.Delete()

$verRange = $tsv.Range("A2:A1048576")
$verRange.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
$verRange.Validation.ErrorTitle = "Value must come from list"
$verRange.Validation.ErrorMessage = "Value must be one of: 1."
$verRange.Validation.IgnoreBlank = $true
$verRange.Validation.InCellDropdown = $true

foreach ($rule in $existing) {
    $rng = $tsv.Range($rule.Col + "2:" + $rule.Col + "1048576")
    if ($rule.Type -eq "decimal") {
        $rng.Validation.Add(2, 1, 1, $rule.Formula1, $rule.Formula2)
    } else {
        $rng.Validation.Add(3, 1, 1, $rule.Formula1)
    }
    $rng.Validation.ErrorTitle = $rule.Title
    $rng.Validation.ErrorMessage = $rule.Msg
    $rng.Validation.IgnoreBlank = $true
    $rng.Validation.InCellDropdown = $true
}

# Keep "Export as TSV" as the active/selected sheet (inserting the new
# worksheet earlier made the runtime switch focus to it).
$tsv.Range("A1").Select()
$tsv.Activate()

Write-Host "Edit applied successfully."
